# fix typo in label: "per capita" -> "per cap." in Predictor column (C)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value  = "ln(GDP [dollars per cap.])"
$ws.Range("C3").Value  = "ln(Tourism - Inbound [per cap.])"
$ws.Range("C4").Value  = "ln(ProMed Mentions [per cap.])"
$ws.Range("C8").Value  = "ln(Publication Bias Index [per cap.])"
$ws.Range("C9").Value  = "ln(AB Exports [dollars per cap.])"
$ws.Range("C11").Value = "ln(Migrant Population [per cap.])"
$ws.Range("C13").Value = "ln(ProMed Mentions [per cap.])"
$ws.Range("C15").Value = "ln(Publication Bias Index [per cap.])"
$ws.Range("C17").Value = "ln(GDP [dollars per cap.])"
